$wb = $excel.ActiveWorkbook

# --- Rename the "Fidelity" sheet to "Investments" ---
$ws = $wb.Worksheets.Item("Fidelity")
$ws.Name = "Investments"

# --- Update the two headline cells that used to carry a "Fidelity " prefix ---
$ws.Range("A1").Value = "Consolidated Year-to-Date Tax Activity"
$ws.Range("A27").Value = "Consolidated Ordinary Dividends and Distributions"

# Writing a fresh value to A27 drops its original quote-prefixed cell style;
# restore it by stamping the (untouched, identically-styled) A47 format back
# onto A27 without touching its newly-written content.
$ws.Range("A47").Copy()
$ws.Range("A27").PasteSpecial(-4122)

# --- Renaming the sheet re-stringifies every OI_ROW defined name and drops
#     both the sheet qualifier on the broken #REF! and the _xlfn. prefix;
#     restore each one explicitly (one at a time, by collection index, so
#     the fix-up of one name doesn't get clobbered by the recompute
#     triggered by the next - looking names up by their "Sheet!Name" text
#     instead of index re-breaks earlier fix-ups for this workbook). ---
$wb.Names.Item(3).RefersTo = "=_xlfn.XMATCH(#REF!,Brackets!#REF!,-1,2)"
$wb.Names.Item(4).RefersTo = "=_xlfn.XMATCH(#REF!,Investments!#REF!,-1,2)"
$wb.Names.Item(5).RefersTo = "=_xlfn.XMATCH(#REF!,Outputsx!#REF!,-1,2)"
$wb.Names.Item(6).RefersTo = "=_xlfn.XMATCH(#REF!,Statics!#REF!,-1,2)"
$wb.Names.Item(7).RefersTo = "=_xlfn.XMATCH(#REF!,#REF!,-1,2)"

# --- Move the active tab from "Statics" to "Investments" and update the
#     remembered selection on that sheet. ---
$ws.Activate()
$ws.Range("B43").Select()
